$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.309.92'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '3.556.83'
$ws.Range('E3').Value = '  +1.83%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '617.57'
$ws.Range('E5').Value = '  +2.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.40'
$ws.Range('E6').Value = '  +4.01%  '
$ws.Range('D7').Value = '3.555.53'
$ws.Range('E7').Value = '  +1.84%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('E9').Value = '  +2.37%  '
$ws.Range('E10').Value = '  +6.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.57'
$ws.Range('E11').Value = '  +8.37%  '
$ws.Range('E12').Value = '  +3.92%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '33.26'
$ws.Range('E13').Value = '  +6.19%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000221'
$ws.Range('E14').Value = '  +2.96%  '
$ws.Range('D15').Value = '4.160.95'
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').Value = '3.566.99'
$ws.Range('E16').Value = '  +2.05%  '
$ws.Range('D17').Value = '68.271.57'
$ws.Range('E17').Value = '  +1.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.117'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.68'
$ws.Range('E19').Value = '  +5.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.95'
$ws.Range('E20').Value = '  +6.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.09'
$ws.Range('E21').Value = '  +12.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '455.19'
$ws.Range('E22').Value = '  +2.43%  '
$ws.Range('E23').Value = '  +4.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.45'
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('E25').Value = '  +5.17%  '
$ws.Range('D26').Value = '3.699.21'
$ws.Range('E26').Value = '  +1.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.25'
$ws.Range('E28').Value = '  +13.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.46'
$ws.Range('E29').Value = '  +4.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.71'
$ws.Range('E30').Value = '  +12.15%  '
$ws.Range('E31').Value = '  +4.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.173'
$ws.Range('E32').Value = '  +5.23%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.35'
$ws.Range('E34').Value = '  +5.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.08'
$ws.Range('E35').Value = '  +1.74%  '
$ws.Range('E36').Value = '  +4.98%  '
$ws.Range('D37').Value = '3.552.76'
$ws.Range('E37').Value = '  +2.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.19'
$ws.Range('E38').Value = '  +2.83%  '
$ws.Range('E39').Value = '  +10.35%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '179.59'
$ws.Range('E41').Value = '  +2.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0920'
$ws.Range('E42').Value = '  +5.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.66'
$ws.Range('E44').Value = '  +5.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '31.41'
$ws.Range('E45').Value = '  +12.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.894'
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.35'
$ws.Range('E47').Value = '  +8.09%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.15'
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.64'
$ws.Range('E49').Value = '  +5.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.78'
$ws.Range('E50').Value = '  +3.57%  '
$ws.Range('E51').Value = '  +8.41%  '
